# Fix ticker typo: APPL -> AAPL in both the "rsu" and "sell_orders" sheets.
$wb = $excel.ActiveWorkbook

$wsRsu = $wb.Worksheets.Item("rsu")
$wsRsu.Range("B6").Value = "AAPL"

$wsSell = $wb.Worksheets.Item("sell_orders")
$wsSell.Range("B6").Value = "AAPL"

# Reproduce the author's final UI state when saving: they had just fixed
# cell B7 on "rsu" (selected it) before finally landing on the
# "money_transfers" tab.
$wsRsu.Activate()
$wsRsu.Range("B7").Select()

$wsMoneyTransfers = $wb.Worksheets.Item("money_transfers")
$wsMoneyTransfers.Activate()
